$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Multiple reason" feature: append a second comma-separated reason to each
# existing value in column D (the "Why did you join this college?" answers).
$ws.Range("D2").Value = "Good Academics,Good Placements"
$ws.Range("D3").Value = "Good Placements,Good Infrastructure or Facilities"
$ws.Range("D4").Value = "Reputation or Brand,Near to Home"
$ws.Range("D5").Value = "Near to Home,Good Placements"
$ws.Range("D6").Value = "Sports Facilities,Reputation or Brand"
$ws.Range("D7").Value = "Good Infrastructure or Facilities,CET"
$ws.Range("D8").Value = "CET,Good Academics"

# Widen column D to fit the longer combined text (OOXML width 24.25 -> 38.0).
# ColumnWidth of 37.1667 chars serializes to an OOXML width of exactly 38.0.
$ws.Columns.Item(4).ColumnWidth = 37.16666666666667
